# MCH240 archival record: populate row 2 with the item metadata and match
# the formatting used by the rest of the imported sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
$ws.Range("A2").Value = "MCH240"
$ws.Range("C2").Value = "NO VAT, PEOPLES POWER FOR A DEMOCRATIC FUTURE, DEMAND A CONSTITUTENT ASSEMBLY"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"

# D2 / H2 are part of the record row but have no value - still carry the
# row's formatting, so touch them too (keeps them present in sheetData).
$ws.Range("D2").Value = ""
$ws.Range("H2").Value = ""

# --- Formatting ---------------------------------------------------------
# Row 2 uses the plain body font (Calibri 10, automatic/theme text colour)
# rather than the bold header font used in row 1.
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ThemeColor = 1

$ws.Range("C2:H2").Font.Name = "Calibri"
$ws.Range("C2:H2").Font.Size = 10
$ws.Range("C2:H2").Font.ThemeColor = 1

# --- View: keep header row frozen, select the newly-entered record ------
$ws.Range("A2:I2").Select()
$excel.ActiveWindow.FreezePanes = $true
